# Add a new row (row 7) to the "反转一个单链表" (reverse a singly linked list)
# table of leetcode linked-list problems, and move the active selection onto
# the newly filled row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new data row. Columns: No. | leetcode | 题目 | 解题方法 | 解题关键词 | 时间复杂度 | 空间复杂度
# The order below matters: it controls the order new entries are appended to
# the shared-string table (解题方法, then 题目, then 解题关键词), matching the
# target workbook's shared string indices.
$ws.Range("D7").Value = "1 三个指针，分别指向当前节点cur，前一个节点prev，下一个节点next`n2 当前节点不为空，就向下执行；为空，就结束循环，返回`n3 保存cur节点的下一个节点到next`n4 当前指针cur指向prev`n5 当前指针cur赋值给prev`n6 下一个节点next赋值给cur指针"
$ws.Range("C7").Value = "反转一个单链表"
$ws.Range("E7").Value = "前后指针`n当前指针`n链表插入/删除"
$ws.Range("F7").Value = "O(N), N是元素个数"
$ws.Range("G7").Value = "O(1)"
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 206

# Row 7 grew tall enough to fit the wrapped solution text.
$ws.Rows.Item(7).RowHeight = 160

# The view now scrolls down a row and the selection moves to the newly
# filled-in solution cell.
$ws.Range("D7").Select()
